$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2412.9167
$ws.Range("I6").Value = 451.1
$ws.Range("J6").Value = 3814.2144
$ws.Range("K6").Value = 1353.3
$ws.Range("L6").Value = 11442.6432
$ws.Range("M6").Value = -1241.3
$ws.Range("N6").Value = -11666.6432

$ws.Range("H17").Value = 434.3846
$ws.Range("J17").Value = 240.63637
$ws.Range("L17").Value = 721.9091100000001
$ws.Range("N17").Value = -1057.90911

$ws.Range("H33").Value = 171.96428
$ws.Range("I33").Value = 154.34616
$ws.Range("K33").Value = 154.34616
$ws.Range("M33").Value = 74.65384

$ws.Range("H51").Value = 6166.6665
$ws.Range("I51").Value = 4666.6665
$ws.Range("J51").Value = 7666.6665
$ws.Range("K51").Value = 4666.6665
$ws.Range("L51").Value = 7666.6665
$ws.Range("M51").Value = -4182.6665
$ws.Range("N51").Value = -8634.666499999999

$ws.Range("H101").Value = 1182.762
$ws.Range("I101").Value = 498.2857
$ws.Range("J101").Value = 1525
$ws.Range("K101").Value = 1494.8571
$ws.Range("L101").Value = 4575
$ws.Range("M101").Value = 127.1428999999998
$ws.Range("N101").Value = -7819

$ws.Range("H138").Value = 2941.5442
$ws.Range("J138").Value = 3031.0317
$ws.Range("L138").Value = 9093.0951
$ws.Range("N138").Value = -19373.0951

$ws.Range("H141").Value = 12328.6
$ws.Range("I141").Value = 16212.286
$ws.Range("K141").Value = 48636.858
$ws.Range("M141").Value = -43456.858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2318.6155
$ws.Range("I45").Value = 2374.0908
$ws.Range("J45").Value = 2013.5
$ws.Range("K45").Value = 2374.0908
$ws.Range("L45").Value = 2013.5
$ws.Range("M45").Value = -1997.0908
$ws.Range("N45").Value = -2767.5

$ws.Range("H61").Value = 1804
$ws.Range("I61").Value = 1804
$ws.Range("K61").Value = 1804
$ws.Range("M61").Value = -1592

$ws.Range("H122").Value = 3464.9788
$ws.Range("I122").Value = 2917.9143
$ws.Range("J122").Value = 5060.5835
$ws.Range("K122").Value = 8753.742899999999
$ws.Range("L122").Value = 15181.7505
$ws.Range("M122").Value = -6303.742899999999
$ws.Range("N122").Value = -20081.7505

$ws.Range("H136").Value = 1804
$ws.Range("I136").Value = 1804
$ws.Range("K136").Value = 5412
$ws.Range("M136").Value = -2862

$ws.Range("H137").Value = 45770
$ws.Range("J137").Value = 45770
$ws.Range("L137").Value = 45770
$ws.Range("N137").Value = -55970

$ws.Range("H138").Value = 76000
$ws.Range("J138").Value = 76000
$ws.Range("L138").Value = 76000
$ws.Range("N138").Value = -86280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1505
$ws.Range("I99").Value = 839.35486
$ws.Range("J99").Value = 3092.3076
$ws.Range("K99").Value = 839.35486
$ws.Range("L99").Value = 3092.3076
$ws.Range("M99").Value = 658.64514
$ws.Range("N99").Value = -6088.3076

$ws.Range("H103").Value = 35324.05
$ws.Range("J103").Value = 35324.05
$ws.Range("L103").Value = 35324.05
$ws.Range("N103").Value = -37668.05

$ws.Range("H137").Value = 49186.668
$ws.Range("J137").Value = 49186.668
$ws.Range("L137").Value = 49186.668
$ws.Range("N137").Value = -59386.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 10103079
$ws.Range("I16").Value = 22223818
$ws.Range("J16").Value = 2462.6667
$ws.Range("K16").Value = 22223818
$ws.Range("L16").Value = 2462.6667
$ws.Range("M16").Value = -22223531
$ws.Range("N16").Value = -3036.6667

$ws.Range("H31").Value = 6915.881
$ws.Range("I31").Value = 3069.8096
$ws.Range("J31").Value = 10761.952
$ws.Range("K31").Value = 3069.8096
$ws.Range("L31").Value = 10761.952
$ws.Range("M31").Value = -2774.8096
$ws.Range("N31").Value = -11351.952

$ws.Range("H34").Value = 6915.881
$ws.Range("I34").Value = 3069.8096
$ws.Range("J34").Value = 10761.952
$ws.Range("K34").Value = 3069.8096
$ws.Range("L34").Value = 10761.952
$ws.Range("M34").Value = -2867.8096
$ws.Range("N34").Value = -11165.952

$ws.Range("H58").Value = 2475.5476
$ws.Range("I58").Value = 1547.5161
$ws.Range("J58").Value = 5090.909
$ws.Range("K58").Value = 1547.5161
$ws.Range("L58").Value = 5090.909
$ws.Range("M58").Value = -1344.5161
$ws.Range("N58").Value = -5496.909

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H94").Value = 1795.1333
$ws.Range("J94").Value = 1948
$ws.Range("L94").Value = 1948
$ws.Range("N94").Value = -2850

$ws.Range("H104").Value = 36866.332
$ws.Range("J104").Value = 36866.332
$ws.Range("L104").Value = 36866.332
$ws.Range("N104").Value = -42108.332

$ws.Range("H113").Value = 10103079
$ws.Range("I113").Value = 22223818
$ws.Range("J113").Value = 2462.6667
$ws.Range("K113").Value = 22223818
$ws.Range("L113").Value = 2462.6667
$ws.Range("M113").Value = -22221648
$ws.Range("N113").Value = -6802.6667

$ws.Range("H119").Value = 33900.2
$ws.Range("J119").Value = 33900.2
$ws.Range("L119").Value = 33900.2
$ws.Range("N119").Value = -43576.2

$ws.Range("H134").Value = 7076.857
$ws.Range("I134").Value = 8685.857
$ws.Range("J134").Value = 3858.8572
$ws.Range("K134").Value = 26057.571
$ws.Range("L134").Value = 11576.5716
$ws.Range("M134").Value = -23522.571
$ws.Range("N134").Value = -16646.5716

$ws.Range("H136").Value = 2475.5476
$ws.Range("I136").Value = 1547.5161
$ws.Range("J136").Value = 5090.909
$ws.Range("K136").Value = 4642.5483
$ws.Range("L136").Value = 15272.727
$ws.Range("M136").Value = -2092.5483
$ws.Range("N136").Value = -20372.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 214.125
$ws.Range("I7").Value = 193.5
$ws.Range("J7").Value = 276
$ws.Range("K7").Value = 580.5
$ws.Range("L7").Value = 828
$ws.Range("M7").Value = -468.5
$ws.Range("N7").Value = -1052

$ws.Range("H92").Value = 297.8
$ws.Range("I92").Value = 297.8
$ws.Range("K92").Value = 893.4000000000001
$ws.Range("M92").Value = 354.5999999999999

$ws.Range("H113").Value = 5000567
$ws.Range("I113").Value = 634.1
$ws.Range("J113").Value = 8333855.5
$ws.Range("K113").Value = 1902.3
$ws.Range("L113").Value = 25001566.5
$ws.Range("M113").Value = 267.6999999999998
$ws.Range("N113").Value = -25005906.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 27816.363
$ws.Range("J46").Value = 28073
$ws.Range("L46").Value = 28073
$ws.Range("N46").Value = -28385

$ws.Range("H96").Value = 29666.666
$ws.Range("J96").Value = 29666.666
$ws.Range("L96").Value = 29666.666
$ws.Range("N96").Value = -35158.666

$ws.Range("H137").Value = 43780
$ws.Range("J137").Value = 43780
$ws.Range("L137").Value = 43780
$ws.Range("N137").Value = -53980

$ws.Range("H141").Value = 38189.855
$ws.Range("J141").Value = 38554.832
$ws.Range("L141").Value = 38554.832
$ws.Range("N141").Value = -48914.832

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5482
$ws.Range("I7").Value = 2076.8
$ws.Range("J7").Value = 7914.2856
$ws.Range("K7").Value = 2076.8
$ws.Range("L7").Value = 7914.2856
$ws.Range("M7").Value = -1964.8
$ws.Range("N7").Value = -8138.2856

$ws.Range("H126").Value = 5482
$ws.Range("I126").Value = 2076.8
$ws.Range("J126").Value = 7914.2856
$ws.Range("K126").Value = 6230.400000000001
$ws.Range("L126").Value = 23742.8568
$ws.Range("M126").Value = -3760.400000000001
$ws.Range("N126").Value = -28682.8568

$ws.Range("H132").Value = 4660.4165
$ws.Range("I132").Value = 3696.9412
$ws.Range("J132").Value = 7000.2856
$ws.Range("K132").Value = 11090.8236
$ws.Range("L132").Value = 21000.8568
$ws.Range("M132").Value = -8560.8236
$ws.Range("N132").Value = -26060.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6680.4707
$ws.Range("I136").Value = 2823.125
$ws.Range("J136").Value = 10109.223
$ws.Range("K136").Value = 8469.375
$ws.Range("L136").Value = 30327.669
$ws.Range("M136").Value = -5919.375
$ws.Range("N136").Value = -35427.669
